$wb = $excel.ActiveWorkbook

$wsApraksts   = $wb.Worksheets.Item("Apraksts")
$wsPrasiba    = $wb.Worksheets.Item("Prasība")
$wsTestpiem   = $wb.Worksheets.Item("Testpiemēri")
$wsZurnals    = $wb.Worksheets.Item("Testēšanas žurnāls")

# ---------------------------------------------------------------------------
# Testpiemēri (sheet3): a new test case (TP.BB.SES.08. / PR.SES.07.) is
# inserted at row 12 - its A/B/H cells carry over the style+data that used
# to live on row 14 (which keeps only its B cell), and D12 gets a brand new
# description string.
# ---------------------------------------------------------------------------

# Preserve formatting of the cells about to move / change by copying styles
# from their current (pre-edit) homes before any values are touched.
$wsTestpiem.Range("A14").Copy()
$wsTestpiem.Range("A12").PasteSpecial(-4122)

$wsTestpiem.Range("H14").Copy()
$wsTestpiem.Range("H12").PasteSpecial(-4122)

$wsTestpiem.Range("B14").Copy()
$wsTestpiem.Range("B12").PasteSpecial(-4122)

$wsTestpiem.Range("A12").Value = "TP.BB.SES.08."
$wsTestpiem.Range("B12").Value = "Reģistrēšanās tīmekļvietnē"
$wsTestpiem.Range("D12").Value = "Reģistrēšanās tīmekļvietnē ievadot ar datubāzes vērtības sakrītošu ierakstu"
$wsTestpiem.Range("H12").Value = "PR.SES.07."
$wsTestpiem.Rows("12").RowHeight = 47.25

# Row 14 loses its A (moved to A12) and H (moved to H12) cells entirely.
$wsTestpiem.Range("A14").Clear()
$wsTestpiem.Range("H14").Clear()

# ---------------------------------------------------------------------------
# Prasība (sheet2): add a new "Administratora panelis" category row at the
# bottom, matching the style of the other category rows (e.g. row 28).
# ---------------------------------------------------------------------------

$wsPrasiba.Range("A28:B28").Copy()
$wsPrasiba.Range("A30:B30").PasteSpecial(-4122)
$wsPrasiba.Range("A30").Value = "Administratora panelis"

# ---------------------------------------------------------------------------
# Testēšanas žurnāls (sheet4): clear the "Veiksmīgs" results for rows 17/18
# (leave the cell / its style in place, only drop the value).
# ---------------------------------------------------------------------------

$wsZurnals.Range("F17").ClearContents()
$wsZurnals.Range("F18").ClearContents()

# ---------------------------------------------------------------------------
# View state: update selections on every sheet, and make "Prasība" the
# active / selected tab last so it ends up as the workbook's active sheet.
# ---------------------------------------------------------------------------

$wsTestpiem.Range("I12").Select() | Out-Null
$wsZurnals.Range("F21").Select() | Out-Null

$wsPrasiba.Activate()
$wsPrasiba.Range("F19").Select() | Out-Null
